# Fruta / hortaliza, semanal
# Insert a new weekly record row above row 45 (shifting existing rows 45-142
# down to 46-143) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 45; everything below shifts down one row.
$ws.Rows("45:45").Insert()

# Populate the new row 45 with the new weekly record.
$ws.Range("A45").Value = 5
$ws.Range("B45").Value = "Macroferia Regional de Talca"
$ws.Range("C45").Value = "Maule"
$ws.Range("D45").Value = 45246
$ws.Range("E45").Value = 7
$ws.Range("F45").Value = "Fruta"
$ws.Range("G45").Value = 100101
$ws.Range("H45").Value = "Berries"
$ws.Range("I45").Value = 100101001
$ws.Range("J45").Value = "Arándano (blue)"
$ws.Range("K45").Value = "Sin especificar"
$ws.Range("L45").Value = "Primera"
$ws.Range("M45").Value = 120
$ws.Range("N45").Value = 6000
$ws.Range("O45").Value = 6000
$ws.Range("P45").Value = 6000
$ws.Range("Q45").Value = "`$/bandeja 2 kilos"
$ws.Range("R45").Value = "Provincia de Curicó"
$ws.Range("S45").Value = 3000
$ws.Range("T45").Value = 2
